$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking price
# strings (e.g. "95.39") are stored as text, matching the source data
# (inline strings), not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").Value = "42.695.15"
$ws.Range("E2").Value = "  -0.32%  "

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").Value = "2.543.57"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("E4").Value = "  +0.04%  "

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").Value = "313.12"
$ws.Range("E5").Value = "  +3.09%  "

# Row 6: 'Solana' -> 'Solana'
$ws.Range("D6").Value = "95.39"
$ws.Range("E6").Value = "  -2.43%  "

# Row 7: 'XRP' -> 'XRP'
$ws.Range("E7").Value = "  +0.45%  "

# Row 8: 'USDC' -> 'USDC'
$ws.Range("E8").Value = "  +0.00%  "

# Row 9: 'Cardano' -> 'Cardano'
$ws.Range("E9").Value = "  -1.79%  "

# Row 10: 'Avalanche' -> 'Avalanche'
$ws.Range("D10").Value = "36.43"
$ws.Range("E10").Value = "  -1.06%  "

# Row 11: 'Dogecoin' -> 'Dogecoin'
$ws.Range("E11").Value = "  -1.74%  "

# Row 12: 'Polkadot' -> 'TRON'
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.114"
$ws.Range("E12").Value = "  -0.38%  "

# Row 13: 'TRON' -> 'Polkadot'
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "7.58"
$ws.Range("E13").Value = "  +0.00%  "

# Row 14: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D14").Value = "2.933.72"
$ws.Range("E14").Value = "  -0.04%  "

# Row 15: 'Chainlink' -> 'Chainlink'
$ws.Range("D15").Value = "15.65"
$ws.Range("E15").Value = "  +3.72%  "

# Row 16: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D16").Value = "2.545.20"
$ws.Range("E16").Value = "  +0.34%  "

# Row 17: 'Polygon' -> 'Polygon'
$ws.Range("D17").Value = "0.863"
$ws.Range("E17").Value = "  -0.34%  "

# Row 18: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D18").Value = "42.696.18"
$ws.Range("E18").Value = "  -0.35%  "

# Row 19: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range("D19").Value = "13.07"
$ws.Range("E19").Value = "  -2.00%  "

# Row 20: 'Uniswap' -> 'Uniswap'
$ws.Range("D20").Value = "6.64"
$ws.Range("E20").Value = "  +1.05%  "

# Row 21: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D21").Value = "0.0₃0970"
$ws.Range("E21").Value = "  -1.86%  "

# Row 22: 'Litecoin' -> 'Litecoin'
$ws.Range("D22").Value = "71.07"

# Row 23: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D23").Value = "256.01"
$ws.Range("E23").Value = "  -0.02%  "

# Row 24: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("D24").Value = "2.96"
$ws.Range("E24").Value = "  +0.36%  "

# Row 25: 'ImmutableX' -> 'ImmutableX'
$ws.Range("E25").Value = "  -1.57%  "

# Row 26: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D26").Value = "27.69"
$ws.Range("E26").Value = "  -1.42%  "

# Row 27: 'Dai' -> 'Dai'
$ws.Range("E27").Value = "  +0.02%  "

# Row 28: 'Toncoin' -> 'Toncoin'
$ws.Range("E28").Value = "  +1.76%  "

# Row 29: 'InjectiveProtocol' -> 'InjectiveProtocol'
$ws.Range("D29").Value = "39.56"
$ws.Range("E29").Value = "  +4.28%  "

# Row 30: 'Cosmos' -> 'Cosmos'
$ws.Range("D30").Value = "10.08"
$ws.Range("E30").Value = "  -0.74%  "

# Row 31: 'Filecoin' -> 'Filecoin'
$ws.Range("D31").Value = "5.96"
$ws.Range("E31").Value = "  -1.96%  "

# Row 32: 'Monero' -> 'Monero'
$ws.Range("D32").Value = "155.89"
$ws.Range("E32").Value = "  -1.06%  "

# Row 33: 'LidoDAOToken' -> 'LidoDAOToken'
$ws.Range("D33").Value = "3.41"
$ws.Range("E33").Value = "  +3.26%  "

# Row 34: 'Celestia' -> 'Celestia'
$ws.Range("D34").Value = "19.48"
$ws.Range("E34").Value = "  +0.02%  "

# Row 35: 'ARBITRUM' -> 'ARBITRUM'
$ws.Range("E35").Value = "  +0.77%  "

# Row 36: 'Hedera' -> 'Hedera'
$ws.Range("D36").Value = "0.0791"
$ws.Range("E36").Value = "  -0.82%  "

# Row 37: 'WEMIXToken' -> 'WEMIXToken'
$ws.Range("E37").Value = "  -0.11%  "

# Row 38: 'Kaspa' -> 'Kaspa'
$ws.Range("E38").Value = "  -3.96%  "

# Row 39: 'EnergySwap' -> 'EnergySwap'
$ws.Range("D39").Value = "24.66"
$ws.Range("E39").Value = "  -4.77%  "

# Row 40: 'Stellar' -> 'Stellar'
$ws.Range("E40").Value = "  -0.35%  "

# Row 41: 'ApeXProtocol' -> 'ApeXProtocol'
$ws.Range("D41").Value = "2.27"
$ws.Range("E41").Value = "  +8.34%  "

# Row 42: 'RenderToken' -> 'NEARProtocol'
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").Value = "3.38"
$ws.Range("E42").Value = "  -0.48%  "

# Row 43: 'NEARProtocol' -> 'RenderToken'
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "3.85"
$ws.Range("E43").Value = "  -0.93%  "

# Row 44: 'VeChain' -> 'VeChain'
$ws.Range("E44").Value = "  -0.62%  "

# Row 45: 'FirstDigitalUSD' -> 'FirstDigitalUSD'
$ws.Range("E45").Value = "  +0.01%  "

# Row 46: 'Maker' -> 'Maker'
$ws.Range("D46").Value = "2.051.54"
$ws.Range("E46").Value = "  -1.90%  "

# Row 47: 'BitcoinSV' -> 'BitcoinSV'
$ws.Range("D47").Value = "85.29"
$ws.Range("E47").Value = "  -2.80%  "

# Row 48: 'FraxShare' -> 'FraxShare'
$ws.Range("D48").Value = "8.93"
$ws.Range("E48").Value = "  +0.80%  "

# Row 49: 'ordi' -> 'ordi'
$ws.Range("D49").Value = "75.44"
$ws.Range("E49").Value = "  +1.31%  "

# Row 50: 'RocketPoolETH' -> 'RocketPoolETH'
$ws.Range("D50").Value = "2.789.79"
$ws.Range("E50").Value = "  -0.07%  "

# Row 51: 'Algorand' -> 'Algorand'
$ws.Range("E51").Value = "  -0.09%  "

# Restore normal (General) styling on column D so no residual text
# formatting is left behind beyond the value changes themselves.
$ws.Range("D2:D51").Style = "Normal"

